$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "taika waititi"
$ws.Range("A4").Select()
